$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (column D) and 1h volume-change percentage
# (column E) figures on the cryptos sheet, per the latest GitHub Actions
# scrape refresh. Column D values are written with NumberFormat forced to
# text ("@") first so Excel does not silently reinterpret price strings
# that look numeric (and lose formatting such as trailing zeros or
# dot-separated thousands groups). Column E percentage strings already
# carry surrounding whitespace, so Excel keeps them as text natively.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.953.45'
$ws.Range("E2").Value = '  -1.01%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.174.74'
$ws.Range("E3").Value = '  -4.40%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '590.41'
$ws.Range("E5").Value = '  -2.65%  '
$ws.Range("E6").Value = '  -4.00%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.174.06'
$ws.Range("E8").Value = '  -4.39%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.513'
$ws.Range("E9").Value = '  -1.32%  '
$ws.Range("E10").Value = '  -4.52%  '
$ws.Range("E11").Value = '  -3.38%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.457'
$ws.Range("E12").Value = '  -2.44%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000239'
$ws.Range("E13").Value = '  -3.63%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.88'
$ws.Range("E14").Value = '  +0.50%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.698.18'
$ws.Range("E15").Value = '  -4.37%  '
$ws.Range("E16").Value = '  -1.95%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.174.90'
$ws.Range("E17").Value = '  -4.34%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '62.912.47'
$ws.Range("E18").Value = '  -1.23%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.62'
$ws.Range("E19").Value = '  -3.48%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '461.64'
$ws.Range("E20").Value = '  -3.88%  '
$ws.Range("E21").Value = '  -1.43%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.713'
$ws.Range("E22").Value = '  -3.03%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.66'
$ws.Range("E23").Value = '  -6.13%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.43'
$ws.Range("E24").Value = '  -2.27%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.48'
$ws.Range("E25").Value = '  -1.82%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.69'
$ws.Range("E27").Value = '  -2.94%  '
$ws.Range("E28").Value = '  +0.09%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.82'
$ws.Range("E29").Value = '  -4.75%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.79'
$ws.Range("E30").Value = '  -6.40%  '
$ws.Range("E31").Value = '  -6.37%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '27.31'
$ws.Range("E32").Value = '  -6.14%  '
$ws.Range("E33").Value = '  -2.83%  '
$ws.Range("E34").Value = '  -6.14%  '
$ws.Range("E35").Value = '  -6.33%  '
$ws.Range("E36").Value = '  -3.09%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '51.10'
$ws.Range("E37").Value = '  -3.04%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0₃0713'
$ws.Range("E38").Value = '  -4.14%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0388'
$ws.Range("E39").Value = '  -2.91%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '406.93'
$ws.Range("E40").Value = '  -6.11%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.70'
$ws.Range("E41").Value = '  -1.88%  '
$ws.Range("E42").Value = '  -2.64%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.113'
$ws.Range("E43").Value = '  -3.29%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.790.62'
$ws.Range("E44").Value = '  -9.81%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.254'
$ws.Range("E45").Value = '  -3.75%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.15'
$ws.Range("E46").Value = '  -2.59%  '
$ws.Range("E47").Value = '  -0.03%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '25.75'
$ws.Range("E48").Value = '  -2.21%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '123.35'
$ws.Range("E49").Value = '  -1.24%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.70'
$ws.Range("E50").Value = '  -5.57%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.112'
$ws.Range("E51").Value = '  -2.30%  '
